# Atualização de bases das ligas, do dia: 28-05-2024 às 19:13
#
# The source rows got re-sorted/re-ordered upstream; as a result the data
# (everything except the running index in column A) for several pairs of
# rows was swapped, and one trio of rows was cyclically rotated.
# This script reproduces that re-shuffle directly on the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($row1, $row2) {
    $range1 = $ws.Range("B$row1`:AD$row1")
    $range2 = $ws.Range("B$row2`:AD$row2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value = $vals2
    $range2.Value = $vals1
}

# Simple pairwise swaps of row content (column A index stays put)
Swap-Rows 9 10
Swap-Rows 50 51
Swap-Rows 88 89
Swap-Rows 92 93
Swap-Rows 107 108
Swap-Rows 119 120

# Three-way rotation: new101 = old102, new102 = old103, new103 = old101
$range101 = $ws.Range("B101:AD101")
$range102 = $ws.Range("B102:AD102")
$range103 = $ws.Range("B103:AD103")

$vals101 = $range101.Value2
$vals102 = $range102.Value2
$vals103 = $range103.Value2

$range101.Value = $vals102
$range102.Value = $vals103
$range103.Value = $vals101
